$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values are stored as text (inline
# strings) in the source workbook, so a leading apostrophe is used to force
# Excel to keep them as text instead of auto-converting to numbers/dates.
$updates = @{
    'D2' = "'261.11"
    'E2' = "'2.04%"
    'D3' = "'27.31"
    'E3' = "'1.20%"
    'D4' = "'4.716"
    'E4' = "'8.98%"
    'E5' = "'3.15%"
    'D6' = "'6.662"
    'E6' = "'0.55%"
    'D7' = "'0.8493"
    'E7' = "'-0.18%"
    'D8' = "'0.9175"
    'E8' = "'-1.82%"
    'D9' = "'0.1402"
    'E9' = "'1.49%"
    'D10' = "'0.04884"
    'E10' = "'2.47%"
    'D11' = "'0.07091"
    'E11' = "'0.13%"
    'D12' = "'0.03145"
    'E12' = "'2.35%"
    'E13' = "'-0.25%"
    'D14' = "'0.001533"
    'E14' = "'0.10%"
    'D15' = "'0.0006087"
    'E15' = "'0.53%"
    'D16' = "'0.006054"
    'E16' = "'-1.62%"
    'E17' = "'-1.02%"
    'E18' = "'-0.56%"
    'E20' = "'2.51%"
    'D21' = "'0.1289"
    'E21' = "'1.44%"
    'D22' = "'4.089"
    'E22' = "'4.46%"
    'E23' = "'-0.38%"
    'D24' = "'0.001219"
    'E24' = "'-0.07%"
    'E25' = "'-8.58%"
    'D26' = "'0.0001200"
    'E26' = "'0.00%"
    'E27' = "'3.02%"
    'D40' = "'0.03870"
    'E40' = "'1.56%"
    'D41' = "'0.1113"
    'E41' = "'1.19%"
    'D42' = "'0.004130"
    'E42' = "'6.49%"
    'D43' = "'0.01635"
    'E43' = "'17.39%"
    'E44' = "'-9.12%"
    'D45' = "'0.00005325"
    'E45' = "'-0.76%"
    'E46' = "'-0.07%"
    'E47' = "'-0.87%"
    'D48' = "'0.1353"
    'E48' = "'-46.24%"
    'E49' = "'-0.07%"
    'E50' = "'-0.07%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
